$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 22) below the existing data (last row was 21),
# mirroring the same "Chirimoya" record shape as the rows above it.
$row = 22

$ws.Range("A$row").Value = 7
$ws.Range("B$row").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C$row").Value = "Ñuble"

$ws.Range("D$row").Value = 45191
$ws.Range("D$row").NumberFormat = $ws.Range("D21").NumberFormat

$ws.Range("E$row").Value = 16
$ws.Range("F$row").Value = "Fruta"
$ws.Range("G$row").Value = 100107
$ws.Range("H$row").Value = "Otros"
$ws.Range("I$row").Value = 100107002
$ws.Range("J$row").Value = "Chirimoya"
$ws.Range("K$row").Value = "Cultivar IV Región"
$ws.Range("L$row").Value = "Primera"
$ws.Range("M$row").Value = 30
$ws.Range("N$row").Value = 21000
$ws.Range("O$row").Value = 21000
$ws.Range("P$row").Value = 21000
$ws.Range("Q$row").Value = "$/bandeja 10 kilos"
$ws.Range("R$row").Value = "Provincia de Limarí"
$ws.Range("S$row").Value = 2100
$ws.Range("T$row").Value = 10
